$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 30   Number  52"
$ws.Range("C9").Value = "Report Covering the Week  12/25/2023  Through  12/31/2023"

# --- Cells that must become text placeholders ("0" / "***.*") ---
# NumberFormat is forced to text ("@") before assignment so the purely-numeric-looking
# string "0" is stored as a shared string rather than re-parsed into a number, then the
# "@"/General format is not needed afterwards: the stored type is already text once written.
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "0"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "***.*"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "***.*"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"

# --- Plain numeric updates ---
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 150
$ws.Range("I15").Value = 30
$ws.Range("K15").Value = -21.052631578947
$ws.Range("L15").Value = 36.363636363636
$ws.Range("M15").Value = 15.384615384615
$ws.Range("N15").Value = -61.538461538461
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -46.153846153846
$ws.Range("I16").Value = 140
$ws.Range("J16").Value = 133
$ws.Range("K16").Value = 5.263157894736
$ws.Range("L16").Value = -13.043478260869
$ws.Range("M16").Value = -62.162162162162
$ws.Range("N16").Value = -88.105352591333
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = -3.333333333333
$ws.Range("I17").Value = 361
$ws.Range("J17").Value = 389
$ws.Range("K17").Value = -7.197943444730
$ws.Range("L17").Value = -9.068010075566
$ws.Range("M17").Value = 1.404494382022
$ws.Range("N17").Value = -52.374670184696
$ws.Range("C18").Value = 4
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 71.428571428571
$ws.Range("I18").Value = 86
$ws.Range("K18").Value = -31.2
$ws.Range("L18").Value = -32.8125
$ws.Range("M18").Value = -79.523809523809
$ws.Range("N18").Value = -92.521739130434
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -62.5
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = -36.111111111111
$ws.Range("I19").Value = 376
$ws.Range("J19").Value = 446
$ws.Range("K19").Value = -15.695067264574
$ws.Range("L19").Value = 1.897018970189
$ws.Range("M19").Value = -45.743145743145
$ws.Range("N19").Value = -90.557508789553
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = 12.5
$ws.Range("I20").Value = 233
$ws.Range("J20").Value = 240
$ws.Range("K20").Value = -2.916666666666
$ws.Range("L20").Value = 16.5
$ws.Range("M20").Value = -11.742424242424
$ws.Range("N20").Value = -87.288597926895
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -29.629629629629
$ws.Range("F21").Value = 105
$ws.Range("G21").Value = 113
$ws.Range("H21").Value = -7.079646017699
$ws.Range("I21").Value = 1234
$ws.Range("J21").Value = 1380
$ws.Range("K21").Value = -10.579710144927
$ws.Range("L21").Value = -4.341085271317
$ws.Range("M21").Value = -42.551210428305
$ws.Range("N21").Value = -86.314738826660
$ws.Range("I23").Value = 13
$ws.Range("K23").Value = 116.666666666667
$ws.Range("L23").Value = -23.529411764705
$ws.Range("M23").Value = -13.333333333333
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -43.333333333333
$ws.Range("F24").Value = 108
$ws.Range("G24").Value = 107
$ws.Range("H24").Value = 0.934579439252
$ws.Range("I24").Value = 1188
$ws.Range("J24").Value = 1329
$ws.Range("K24").Value = -10.609480812641
$ws.Range("L24").Value = 19.758064516129
$ws.Range("M24").Value = 7.803992740471
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 10
$ws.Range("F25").Value = 58
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = 52.631578947368
$ws.Range("I25").Value = 696
$ws.Range("J25").Value = 539
$ws.Range("K25").Value = 29.128014842300
$ws.Range("L25").Value = 43.209876543209
$ws.Range("M25").Value = -12.67252195734
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 150
$ws.Range("I26").Value = 48
$ws.Range("K26").Value = -9.433962264150
$ws.Range("L26").Value = 17.073170731707
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 33.333333333333
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 8
$ws.Range("I28").Value = 27
$ws.Range("K28").Value = -12.903225806451
$ws.Range("L28").Value = -52.631578947368
$ws.Range("M28").Value = -55
$ws.Range("N28").Value = -83.125
$ws.Range("F29").Value = 4
$ws.Range("I29").Value = 18
$ws.Range("K29").Value = -28
$ws.Range("L29").Value = -59.090909090909
$ws.Range("M29").Value = -64
$ws.Range("N29").Value = -87.5
